$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.559.99'
$ws.Range('E2').Value = '  +4.53%  '
$ws.Range('D3').Value = '3.486.91'
$ws.Range('E3').Value = '  +2.81%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '589.39'
$ws.Range('E5').Value = '  +3.33%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '168.68'
$ws.Range('E6').Value = '  +4.52%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '3.486.61'
$ws.Range('E8').Value = '  +2.80%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.592'
$ws.Range('E9').Value = '  +7.77%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.32'
$ws.Range('E10').Value = '  +0.69%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.126'
$ws.Range('E11').Value = '  +6.37%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.438'
$ws.Range('E12').Value = '  +3.91%  '
$ws.Range('D13').Value = '4.093.96'
$ws.Range('E13').Value = '  +3.14%  '
$ws.Range('E14').Value = '  -0.13%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '28.15'
$ws.Range('E15').Value = '  +4.50%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000178'
$ws.Range('E16').Value = '  +3.52%  '
$ws.Range('D17').Value = '66.585.32'
$ws.Range('E17').Value = '  +4.51%  '
$ws.Range('D18').Value = '3.485.36'
$ws.Range('E18').Value = '  +3.48%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.33'
$ws.Range('E19').Value = '  +3.64%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.08'
$ws.Range('E20').Value = '  +3.75%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '391.99'
$ws.Range('E21').Value = '  +4.36%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.92'
$ws.Range('E22').Value = '  +2.07%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '73.22'
$ws.Range('E23').Value = '  +4.33%  '
$ws.Range('E24').Value = '  -0.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.536'
$ws.Range('E25').Value = '  +4.34%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000121'
$ws.Range('E26').Value = '  +5.69%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.09'
$ws.Range('E27').Value = '  +6.89%  '
$ws.Range('E28').Value = '  +2.17%  '
$ws.Range('E29').Value = '  -0.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.34'
$ws.Range('E30').Value = '  +5.03%  '
$ws.Range('E31').Value = '  +6.09%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.06'
$ws.Range('E32').Value = '  +3.02%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '23.60'
$ws.Range('E33').Value = '  +3.37%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.41'
$ws.Range('E34').Value = '  +4.98%  '
$ws.Range('E35').Value = '  +0.10%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.62'
$ws.Range('E36').Value = '  +9.89%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '161.91'
$ws.Range('E37').Value = '  +1.59%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.888'
$ws.Range('E38').Value = '  +3.44%  '
$ws.Range('E39').Value = '  +6.37%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.66'
$ws.Range('E40').Value = '  +6.72%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0745'
$ws.Range('E41').Value = '  +3.26%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.72'
$ws.Range('E42').Value = '  +4.76%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '26.49'
$ws.Range('E43').Value = '  +3.76%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '27.10'
$ws.Range('E44').Value = '  +4.35%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '43.20'
$ws.Range('E45').Value = '  +1.08%  '
$ws.Range('D46').Value = '2.781.10'
$ws.Range('E46').Value = '  +1.42%  '
$ws.Range('E47').Value = '  +2.27%  '
$ws.Range('B48').Value = 'Bittensor'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '348.50'
$ws.Range('E48').Value = '  +6.64%  '
$ws.Range('B49').Value = 'dogwifhat'
$ws.Range('C49').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.47'
$ws.Range('E49').Value = '  +2.42%  '
$ws.Range('E50').Value = '  +5.95%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '33.93'
$ws.Range('E51').Value = '  +14.17%  '
